$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4179.55
$ws.Range("I15").Value = 4179.55
$ws.Range("K15").Value = 12538.65
$ws.Range("M15").Value = -12369.65
$ws.Range("H129").Value = 878.8049
$ws.Range("I129").Value = 697.2
$ws.Range("J129").Value = 904.0278
$ws.Range("K129").Value = 2091.6
$ws.Range("L129").Value = 2712.0834
$ws.Range("M129").Value = 2908.4
$ws.Range("N129").Value = -12712.0834
$ws.Range("H138").Value = 3313.25
$ws.Range("I138").Value = 3487.2778
$ws.Range("K138").Value = 10461.8334
$ws.Range("M138").Value = -5321.8334
$ws.Range("H141").Value = 1274812
$ws.Range("I141").Value = 1648109.9
$ws.Range("K141").Value = 4944329.699999999
$ws.Range("M141").Value = -4939149.699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3113.4915
$ws.Range("I32").Value = 2909.9788
$ws.Range("J32").Value = 3910.5833
$ws.Range("K32").Value = 2909.9788
$ws.Range("L32").Value = 3910.5833
$ws.Range("M32").Value = -2622.9788
$ws.Range("N32").Value = -4484.5833
$ws.Range("H45").Value = 1597.6364
$ws.Range("I45").Value = 1429.6666
$ws.Range("K45").Value = 1429.6666
$ws.Range("M45").Value = -1052.6666
$ws.Range("H61").Value = 2647.0571
$ws.Range("I61").Value = 1860.5927
$ws.Range("J61").Value = 5301.375
$ws.Range("K61").Value = 1860.5927
$ws.Range("L61").Value = 5301.375
$ws.Range("M61").Value = -1648.5927
$ws.Range("N61").Value = -5725.375
$ws.Range("H102").Value = 1320
$ws.Range("I102").Value = 981.8
$ws.Range("K102").Value = 981.8
$ws.Range("M102").Value = 640.2
$ws.Range("H136").Value = 2647.0571
$ws.Range("I136").Value = 1860.5927
$ws.Range("J136").Value = 5301.375
$ws.Range("K136").Value = 5581.7781
$ws.Range("L136").Value = 15904.125
$ws.Range("M136").Value = -3031.7781
$ws.Range("N136").Value = -21004.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4940.884
$ws.Range("I134").Value = 6008.533
$ws.Range("J134").Value = 2477.077
$ws.Range("K134").Value = 18025.599
$ws.Range("L134").Value = 7431.231000000001
$ws.Range("M134").Value = -15490.599
$ws.Range("N134").Value = -12501.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2505.6667
$ws.Range("I31").Value = 3103
$ws.Range("K31").Value = 3103
$ws.Range("M31").Value = -2808
$ws.Range("H34").Value = 2505.6667
$ws.Range("I34").Value = 3103
$ws.Range("K34").Value = 3103
$ws.Range("M34").Value = -2901
$ws.Range("H58").Value = 1813343.4
$ws.Range("I58").Value = 2718958
$ws.Range("J58").Value = 2114.25
$ws.Range("K58").Value = 2718958
$ws.Range("L58").Value = 2114.25
$ws.Range("M58").Value = -2718755
$ws.Range("N58").Value = -2520.25
$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 1650
$ws.Range("J99").Value = 2350
$ws.Range("K99").Value = 1650
$ws.Range("L99").Value = 2350
$ws.Range("M99").Value = -152
$ws.Range("N99").Value = -5346
$ws.Range("H107").Value = 1617.7273
$ws.Range("I107").Value = 1756.2222
$ws.Range("K107").Value = 1756.2222
$ws.Range("M107").Value = 163.7778000000001
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 1650
$ws.Range("J126").Value = 2350
$ws.Range("K126").Value = 4950
$ws.Range("L126").Value = 7050
$ws.Range("M126").Value = -2480
$ws.Range("N126").Value = -11990
$ws.Range("H132").Value = 1187.4222
$ws.Range("I132").Value = 832.9737
$ws.Range("J132").Value = 3111.5715
$ws.Range("K132").Value = 2498.9211
$ws.Range("L132").Value = 9334.7145
$ws.Range("M132").Value = 31.07889999999998
$ws.Range("N132").Value = -14394.7145
$ws.Range("H134").Value = 1648.3024
$ws.Range("I134").Value = 1434.3549
$ws.Range("K134").Value = 4303.0647
$ws.Range("M134").Value = -1768.0647
$ws.Range("H136").Value = 1813343.4
$ws.Range("I136").Value = 2718958
$ws.Range("J136").Value = 2114.25
$ws.Range("K136").Value = 8156874
$ws.Range("L136").Value = 6342.75
$ws.Range("M136").Value = -8154324
$ws.Range("N136").Value = -11442.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 2000
$ws.Range("K80").Value = 6000
$ws.Range("M80").Value = -5064
$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 2000
$ws.Range("K83").Value = 18000
$ws.Range("M83").Value = -13320
$ws.Range("H113").Value = 6698.1763
$ws.Range("J113").Value = 913
$ws.Range("L113").Value = 2739
$ws.Range("N113").Value = -7079
$ws.Range("H131").Value = 22374.945
$ws.Range("J131").Value = 27506.834
$ws.Range("L131").Value = 82520.50199999999
$ws.Range("N131").Value = -92600.50199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1041.2727
$ws.Range("I113").Value = 926.6667
$ws.Range("J113").Value = 1084.25
$ws.Range("K113").Value = 926.6667
$ws.Range("L113").Value = 1084.25
$ws.Range("M113").Value = 1243.3333
$ws.Range("N113").Value = -5424.25
$ws.Range("H132").Value = 1133884.2
$ws.Range("I132").Value = 2026133.9
$ws.Range("J132").Value = 3701.4
$ws.Range("K132").Value = 6078401.699999999
$ws.Range("L132").Value = 11104.2
$ws.Range("M132").Value = -6075871.699999999
$ws.Range("N132").Value = -16164.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1461.6154
$ws.Range("I100").Value = 1127.3636
$ws.Range("K100").Value = 1127.3636
$ws.Range("M100").Value = -586.3635999999999
$ws.Range("H136").Value = 2973.276
$ws.Range("I136").Value = 1817.45
$ws.Range("J136").Value = 5541.778
$ws.Range("K136").Value = 5452.35
$ws.Range("L136").Value = 16625.334
$ws.Range("M136").Value = -2902.35
$ws.Range("N136").Value = -21725.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 30000
$ws.Range("J70").Value = 30000
$ws.Range("L70").Value = 30000
$ws.Range("N70").Value = -30630
$ws.Range("H73").Value = 30000
$ws.Range("J73").Value = 30000
$ws.Range("L73").Value = 30000
$ws.Range("N73").Value = -32184
$ws.Range("H132").Value = 1785.3529
$ws.Range("I132").Value = 1042.619
$ws.Range("J132").Value = 2985.1538
$ws.Range("K132").Value = 3127.857
$ws.Range("L132").Value = 8955.4614
$ws.Range("M132").Value = -597.857
$ws.Range("N132").Value = -14015.4614
